$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 15:22"

# --- Helper: write a full data row (country name + 7 numeric columns) ---
function Set-Row($r, $name, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

# Row 26: Noruega - refreshed numbers
Set-Row 26 "Noruega" 6233 14 32 6093 70 0 108

# Rows 43-44: "Serbia" moves above "Emiratos Arabes Unidos" in the country
# list; Serbia gets freshly updated figures while Emiratos Arabes Unidos
# keeps its previous figures but shifts down one row.
Set-Row 43 "Serbia" 3105 238 118 2916 136 5 71
Set-Row 44 "Emiratos Arabes Unidos" 2990 0 268 2708 1 0 14

# Row 73: Barein - refreshed numbers
Set-Row 73 "Barein" 913 26 530 377 3 1 6

# Row 118: Kenia - refreshed numbers
Set-Row 118 "Kenia" 189 5 22 160 2 0 7

# Rows 122-123: "Republica de Yibuti" moves above "Guadalupe"; Yibuti gets
# freshly updated figures while Guadalupe keeps its previous figures but
# shifts down one row.
Set-Row 122 "Republica de Yibuti" 150 10 36 113 0 0 1
Set-Row 123 "Guadalupe" 143 0 67 68 13 0 8

# Row 139: Barbados - refreshed numbers
Set-Row 139 "Barbados" 66 0 11 51 4 1 4

# Rows 190-193: "Belice" moves above "Surinam" and gets freshly updated
# figures; "Surinam" shifts down one row keeping its own previous
# figures; "Malaui" and "Nepal" swap places with each other, each
# keeping its own previous figures.
Set-Row 190 "Belice" 10 1 0 9 1 0 1
Set-Row 191 "Surinam" 10 0 4 5 0 0 1
Set-Row 192 "Malaui" 9 1 0 8 1 0 1
Set-Row 193 "Nepal" 9 0 1 8 0 0 0

# Rows 208-211: "Sudan del Sur" moves above "Anguila"; all of these rows
# happen to carry identical figures (3,0,0,3,0,0,0) so only the country
# names shift while the numbers stay the same.
Set-Row 208 "Sudan del Sur" 3 0 0 3 0 0 0
Set-Row 209 "Anguila" 3 0 0 3 0 0 0
Set-Row 210 "Burundi" 3 0 0 3 0 0 0
Set-Row 211 "Islas Virgenes Britanicas" 3 0 0 3 0 0 0

# Rows 215-216: "San Pedro y Miquelon" swaps above "Yemen", each keeping
# its own previous figures.
Set-Row 215 "San Pedro y Miquelon" 1 0 0 1 0 0 0
Set-Row 216 "Yemen" 1 1 0 1 0 0 0
